$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on Price (D) cells whose new values are plain numeric
# strings, so Excel does not auto-convert them to floating point numbers
# (which would lose exact text representation such as trailing zeros).
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = '26.253.95'
$ws.Range("E2").Value = '  +0.26%  '
$ws.Range("D3").Value = '1.679.80'
$ws.Range("E3").Value = '  +0.66%  '
$ws.Range("E4").Value = '  +0.31%  '
$ws.Range("D5").Value = '218.53'
$ws.Range("E5").Value = '  +0.69%  '
$ws.Range("D6").Value = '0.5277'
$ws.Range("E6").Value = '  +2.99%  '
$ws.Range("E7").Value = '  +0.29%  '
$ws.Range("D8").Value = '0.2706'
$ws.Range("E8").Value = '  +2.40%  '
$ws.Range("D9").Value = '0.06490'
$ws.Range("E9").Value = '  +1.34%  '
$ws.Range("D10").Value = '21.96'
$ws.Range("E10").Value = '  +1.57%  '
$ws.Range("D11").Value = '0.07524'
$ws.Range("E11").Value = '  +1.46%  '
$ws.Range("D12").Value = '1.693.68'
$ws.Range("E12").Value = '  +1.17%  '
$ws.Range("D13").Value = '4.527'
$ws.Range("E13").Value = '  +0.22%  '
$ws.Range("D14").Value = '0.5810'
$ws.Range("E14").Value = '  -0.22%  '
$ws.Range("D15").Value = '0.000008514'
$ws.Range("E15").Value = '  -0.89%  '
$ws.Range("D16").Value = '64.65'
$ws.Range("E16").Value = '  +0.57%  '
$ws.Range("D17").Value = '26.310.11'
$ws.Range("E17").Value = '  +0.33%  '
$ws.Range("D18").Value = '4.927'
$ws.Range("E18").Value = '  -0.30%  '
$ws.Range("E19").Value = '  +0.24%  '
$ws.Range("D20").Value = '10.88'
$ws.Range("E20").Value = '  +0.43%  '
$ws.Range("D21").Value = '190.38'
$ws.Range("E21").Value = '  -0.20%  '
$ws.Range("D22").Value = '6.201'
$ws.Range("E22").Value = '  -0.16%  '
$ws.Range("E23").Value = '  +0.24%  '
$ws.Range("D24").Value = '145.48'
$ws.Range("E24").Value = '  +0.05%  '
$ws.Range("D25").Value = '7.810'
$ws.Range("E25").Value = '  +2.35%  '
$ws.Range("E26").Value = '  +4.31%  '
$ws.Range("E27").Value = '  +1.02%  '
$ws.Range("D28").Value = '0.06567'
$ws.Range("E28").Value = '  +2.48%  '
$ws.Range("E29").Value = '  +4.31%  '
$ws.Range("D30").Value = '1.331'
$ws.Range("E30").Value = '  +0.92%  '
$ws.Range("D31").Value = '3.601'
$ws.Range("D32").Value = '3.592'
$ws.Range("E32").Value = '  +1.94%  '
$ws.Range("E33").Value = '  +1.16%  '
$ws.Range("E34").Value = '  +1.70%  '
$ws.Range("D35").Value = '0.6230'
$ws.Range("E35").Value = '  +2.26%  '
$ws.Range("E36").Value = '  +1.52%  '
$ws.Range("D37").Value = '2.734'
$ws.Range("E37").Value = '  +2.95%  '
$ws.Range("D38").Value = '6.451'
$ws.Range("E38").Value = '  +4.76%  '
$ws.Range("D39").Value = '1.112.75'
$ws.Range("E39").Value = '  +2.67%  '
$ws.Range("E40").Value = '  +1.14%  '
$ws.Range("D41").Value = '0.8762'
$ws.Range("E41").Value = '  +1.29%  '
$ws.Range("D42").Value = '1.015'
$ws.Range("E42").Value = '  +0.51%  '
$ws.Range("D43").Value = '100.81'
$ws.Range("E43").Value = '  -0.45%  '
$ws.Range("D44").Value = '1.830.73'
$ws.Range("E44").Value = '  +0.79%  '
$ws.Range("D45").Value = '0.00000000110'
$ws.Range("E45").Value = '  -1.18%  '
$ws.Range("D46").Value = '57.04'
$ws.Range("E46").Value = '  +1.33%  '
$ws.Range("D47").Value = '8.205'
$ws.Range("E47").Value = '  +1.44%  '
$ws.Range("E48").Value = '  -0.16%  '
$ws.Range("D49").Value = '0.05273'
$ws.Range("E49").Value = '  +1.33%  '
$ws.Range("D50").Value = '6.092'
$ws.Range("E50").Value = '  +3.16%  '
$ws.Range("D51").Value = '0.4292'
$ws.Range("E51").Value = '  +0.00%  '
